$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains one additional forecast-vector row (2007 -> 2008),
# so every existing data row shifts down by one and a new row 53 is added.
# Rewrite the full data block (rows 2-53) with the corrected / updated values.

$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 11.13090654781819
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 9.396507498425466
$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 7.193183327378438
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 9.591339540850875
$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 2008
$ws.Cells.Item(4, 3).Value = 4.672550446571067
$ws.Cells.Item(4, 4).Value = 2009
$ws.Cells.Item(4, 5).Value = 4.422525088127283
$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = -14.96173956806345
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = -4.932343798304595
$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2009
$ws.Cells.Item(6, 3).Value = -14.45332333832743
$ws.Cells.Item(6, 4).Value = 2010
$ws.Cells.Item(6, 5).Value = -2.928447329610073
$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 2.682935444832424
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = -2.225127715916653
$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 2010
$ws.Cells.Item(8, 3).Value = 8.600536527919633
$ws.Cells.Item(8, 4).Value = 2011
$ws.Cells.Item(8, 5).Value = 6.303897256856628
$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 9.399485634179229
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = 1.811802132286955
$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 2011
$ws.Cells.Item(10, 3).Value = 10.25770250047622
$ws.Cells.Item(10, 4).Value = 2012
$ws.Cells.Item(10, 5).Value = 10.22374275635105
$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 5.169490031659674
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 9.213376886330305
$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 2012
$ws.Cells.Item(12, 3).Value = 4.639893381363169
$ws.Cells.Item(12, 4).Value = 2013
$ws.Cells.Item(12, 5).Value = 8.174613408931286
$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = -0.3722371047999662
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 2.684220738731935
$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 2013
$ws.Cells.Item(14, 3).Value = 0.3058963467304165
$ws.Cells.Item(14, 4).Value = 2014
$ws.Cells.Item(14, 5).Value = 2.429116709932622
$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 4.098801479368341
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 2.548306621254004
$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 2014
$ws.Cells.Item(16, 3).Value = 4.068173739091874
$ws.Cells.Item(16, 4).Value = 2015
$ws.Cells.Item(16, 5).Value = 3.9413000500929
$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = 3.75051862559701
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 2.714258593289975
$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 4.984288257750213
$ws.Cells.Item(18, 4).Value = 2016
$ws.Cells.Item(18, 5).Value = 4.188839638544284
$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 2.352205130086071
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 3.873414041014778
$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 2016
$ws.Cells.Item(20, 3).Value = 1.878184267712912
$ws.Cells.Item(20, 4).Value = 2017
$ws.Cells.Item(20, 5).Value = 2.514670279852349
$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 4.083548352538369
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 3.586256146074462
$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = 4.695933104194339
$ws.Cells.Item(22, 4).Value = 2018
$ws.Cells.Item(22, 5).Value = 4.5579527192392
$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 6.704509587264518
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 4.268691600002228
$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 4.861590900330692
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 3.297472770389764
$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 5.402237127943743
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 4.104053120889217
$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 2018
$ws.Cells.Item(26, 3).Value = 4.892602738886098
$ws.Cells.Item(26, 4).Value = 2019
$ws.Cells.Item(26, 5).Value = 1.957202207503861
$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 1.675184815837505
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 3.589879698956056
$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 1.787861866846807
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 4.088367525047842
$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 0.8513583007189629
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 2.225279621195808
$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = 0.8049382522247184
$ws.Cells.Item(30, 4).Value = 2020
$ws.Cells.Item(30, 5).Value = 2.267257846564918
$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 1.015697339178034
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 2.122104735451602
$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = -2.21482332957591
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = -0.6322362079330346
$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -9.810777096850787
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = -6.212835522792448
$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = -8.784173899737169
$ws.Cells.Item(34, 4).Value = 2021
$ws.Cells.Item(34, 5).Value = 2.199380357735481
$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = 5.72229384158125
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = 1.117941783921328
$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = 6.09521976277807
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = 1.839905110456375
$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = 5.797134106720514
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = 2.056896997569879
$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 2021
$ws.Cells.Item(38, 3).Value = 5.110501195359984
$ws.Cells.Item(38, 4).Value = 2022
$ws.Cells.Item(38, 5).Value = 0.3515918738370427
$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 4.526365501075413
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = 0.9036269924846962
$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 3.616930127707629
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = 1.391416039405691
$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 4.232564748995715
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = 2.135688430332006
$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 2022
$ws.Cells.Item(42, 3).Value = 5.120680133083599
$ws.Cells.Item(42, 4).Value = 2023
$ws.Cells.Item(42, 5).Value = 5.934275247805543
$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = 1.102138938525221
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = 4.283383641765459
$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = 0.7171092762090492
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = 2.755142438739822
$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = 0.08070151925247959
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = 1.265176565876436
$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 2023
$ws.Cells.Item(46, 3).Value = -0.5532735011319234
$ws.Cells.Item(46, 4).Value = 2024
$ws.Cells.Item(46, 5).Value = -1.846917864698006
$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = -2.696492768996317
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = 0.4317200868126703
$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = -0.1521036778360019
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 1.645968204809645
$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = -0.9685570952743805
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = -0.01788907424267183
$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 2024
$ws.Cells.Item(50, 3).Value = -1.069674659641462
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = -0.7986414110784379
$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = -4.127040013406502
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = -2.535768303458463
$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = -2.051528019634985
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = -0.3224191428759626
$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 2025
$ws.Cells.Item(53, 3).Value = -2.436529450546909
$ws.Cells.Item(53, 4).Value = 2026
$ws.Cells.Item(53, 5).Value = -0.469872647443903

# Row 53 is brand new; give its date cell (column A) the same number format
# style used by the rest of column A (style index carried via format-only paste).
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
